# Fixed update to excel issue
#
# 1. Rename the "Requested quantity" header on the "Weekly Quantity" sheet
#    to "Weekly_PO_Qty".
# 2. Rename the "Requested quantity" header on the "Monthly Trend" sheet
#    to "Monthly_PO_Qty".
# 3. Add a new "PO Forecast" worksheet (after "Monthly Trend") containing
#    the forecast data (ds / PO_Forecast / yhat_lower / yhat_upper).

$wb = $excel.ActiveWorkbook

# --- 1. Weekly Quantity header -------------------------------------------
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

# --- 2. Monthly Trend header ----------------------------------------------
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 3. New "PO Forecast" sheet, inserted after the last existing sheet ---
$lastIndex = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($lastIndex)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Match the outline + page-margin settings used on the other sheets.
$wsForecast.Outline.SummaryRow = 1
$wsForecast.Outline.SummaryColumn = 1
$wsForecast.PageSetup.LeftMargin = 54
$wsForecast.PageSetup.RightMargin = 54
$wsForecast.PageSetup.TopMargin = 72
$wsForecast.PageSetup.BottomMargin = 72
$wsForecast.PageSetup.HeaderMargin = 36
$wsForecast.PageSetup.FooterMargin = 36

# Header row, styled to match the bold/bordered/centered header used on
# the other sheets.
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

$hdr = $wsForecast.Range("A1:D1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160
$hdr.Borders.LineStyle = 1

# Forecast data rows.
$rows = @(
  @(45123.99999999999, 18, 3.765440548462984, 31.37763204675185),
  @(45130.99999999999, 19, 5.065697155751526, 31.80627926356178),
  @(45144.99999999999, 21, 6.40586444020375, 34.55986708515742),
  @(45270.99999999999, 39, 24.81760960720618, 52.3763705671942),
  @(45277.99999999999, 40, 26.73265582287967, 53.76281252069803),
  @(45298.99999999999, 43, 29.46317537805147, 56.63745430168877),
  @(45305.99999999999, 44, 29.65641698928991, 57.42595593326747),
  @(45312.99999999999, 45, 30.19594602527467, 58.08027871906436),
  @(45319.99999999999, 46, 32.5070705122893, 59.95017303887991),
  @(45326.99999999999, 47, 33.67761552783364, 60.58614570253704),
  @(45333.99999999999, 48, 34.14491679907012, 61.25683669365282),
  @(45340.99999999999, 49, 35.90300242609539, 62.69756896406615),
  @(45347.99999999999, 50, 35.86245854021421, 63.69802792132554),
  @(45354.99999999999, 51, 38.7271123569475, 65.94141655922803)
)

$r = 2
foreach ($row in $rows) {
    $wsForecast.Cells.Item($r, 1).Value = $row[0]
    $wsForecast.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $wsForecast.Cells.Item($r, 2).Value = $row[1]
    $wsForecast.Cells.Item($r, 3).Value = $row[2]
    $wsForecast.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

$wsWeekly.Select()
